$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44363
$ws.Range("J3").Value = 140
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 14500
$ws.Range("P3").Value = 806
$ws.Range("D4").Value = 44412
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 17500
$ws.Range("P4").Value = 972
$ws.Range("D5").Value = 44398
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 17000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 17500
$ws.Range("P5").Value = 972
$ws.Range("D6").Value = 44398
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 15500
$ws.Range("P6").Value = 861
$ws.Range("I7").Value = 'Segunda'
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 17000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 17500
$ws.Range("P7").Value = 972
$ws.Range("D8").Value = 44435
$ws.Range("I8").Value = 'Tercera'
$ws.Range("J8").Value = 120
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 14500
$ws.Range("P8").Value = 806
$ws.Range("D9").Value = 44405
$ws.Range("J9").Value = 140
$ws.Range("I10").Value = 'Segunda'
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 17000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 17500
$ws.Range("P10").Value = 972
$ws.Range("D11").Value = 44433
$ws.Range("H11").Value = 'Cultivar IV Región'
$ws.Range("I11").Value = 'Tercera'
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14500
$ws.Range("N11").Value = '$/bandeja 18 kilos'
$ws.Range("O11").Value = 'Provincia de Limarí'
$ws.Range("P11").Value = 806
$ws.Range("Q11").Value = 18
$ws.Range("D12").Value = 44211
$ws.Range("J12").Value = 140
$ws.Range("K12").Value = 4500
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = 4750
$ws.Range("P12").Value = 475
$ws.Range("D13").Value = 44391
$ws.Range("H13").Value = 'Cultivar IV Región'
$ws.Range("I13").Value = 'Segunda'
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 15500
$ws.Range("N13").Value = '$/bandeja 18 kilos'
$ws.Range("O13").Value = 'Provincia de Limarí'
$ws.Range("P13").Value = 861
$ws.Range("Q13").Value = 18
$ws.Range("D14").Value = 44454
$ws.Range("H14").Value = 'Cultivar IV Región'
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 160
$ws.Range("K14").Value = 19000
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = 19500
$ws.Range("N14").Value = '$/bandeja 18 kilos'
$ws.Range("O14").Value = 'Provincia de Limarí'
$ws.Range("P14").Value = 1083
$ws.Range("Q14").Value = 18
$ws.Range("D15").Value = 44526
$ws.Range("H15").Value = 'Cultivar XV región'
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 5000
$ws.Range("L15").Value = 5500
$ws.Range("M15").Value = 5250
$ws.Range("N15").Value = '$/caja 10 kilos'
$ws.Range("O15").Value = 'Región de Arica y Parinacota'
$ws.Range("P15").Value = 525
$ws.Range("Q15").Value = 10
$ws.Range("D16").Value = 44526
$ws.Range("I16").Value = 'Segunda'
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 4000
$ws.Range("L16").Value = 4500
$ws.Range("M16").Value = 4250
$ws.Range("P16").Value = 425
$ws.Range("D17").Value = 44526
$ws.Range("I17").Value = 'Tercera'
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 3500
$ws.Range("M17").Value = 3250
$ws.Range("P17").Value = 325
$ws.Range("D18").Value = 44533
$ws.Range("H18").Value = 'Cultivar XV región'
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 6000
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 6500
$ws.Range("N18").Value = '$/caja 10 kilos'
$ws.Range("O18").Value = 'Región de Arica y Parinacota'
$ws.Range("P18").Value = 650
$ws.Range("Q18").Value = 10
$ws.Range("D19").Value = 44533
$ws.Range("H19").Value = 'Cultivar XV región'
$ws.Range("I19").Value = 'Segunda'
$ws.Range("J19").Value = 120
$ws.Range("K19").Value = 4000
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = 4500
$ws.Range("N19").Value = '$/caja 10 kilos'
$ws.Range("O19").Value = 'Región de Arica y Parinacota'
$ws.Range("P19").Value = 450
$ws.Range("Q19").Value = 10
$ws.Range("D20").Value = 44221
$ws.Range("J20").Value = 140
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = 5500
$ws.Range("P20").Value = 550
$ws.Range("D21").Value = 44554
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 6000
$ws.Range("M21").Value = 5500
$ws.Range("P21").Value = 550
